# Apply the "new version with timestamp" update to the DaySale report.
#
# Three product rows get updated current-balance / selling-price-total /
# transaction-count figures (reflecting one additional unit sold for each),
# the grand-total cell is bumped to match, and the generated-at timestamp
# in the footer moves forward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # The "selling price total" cells (column P) carry a numeric display
    # format ("0.00") even though the stored cell content is plain text
    # (e.g. "88.0000"). Writing a numeric-looking string straight into
    # .Value would get auto-coerced into a real number (dropping the
    # trailing zeros) and would reset the cell style, so toggle the
    # format to Text for the write and then restore the original format.
    $fmt = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $fmt
}

# Row 15 - CONGESTAL SYRUP 120 ML
#   current balance 2:0 -> 3:0, selling price total 44.0000 -> 88.0000,
#   transactions 1:0 -> 2:0
$ws.Range("H15").Value = "3:0"
Set-TextValue $ws.Range("P15") "88.0000"
$ws.Range("Q15").Value = "2:0"

# Row 39 - NEW-CLAV EXTRA STRENGTH 642.9MG/5ML 50ML SUSP.
#   current balance 1:0 -> 2:0, selling price total 180.0000 -> 270.0000
$ws.Range("H39").Value = "2:0"
Set-TextValue $ws.Range("P39") "270.0000"

# Row 48 - WATER FOR INJECTION AMP. 5 ML
#   current balance 8519:0 -> 8520:0, selling price total 10.0000 -> 12.0000,
#   transactions 5:0 -> 6:0
$ws.Range("H48").Value = "8520:0"
Set-TextValue $ws.Range("P48") "12.0000"
$ws.Range("Q48").Value = "6:0"

# Grand total row picks up the combined +136 delta from the three rows above
$ws.Range("P61").Value = 2765.495

# Footer generation timestamp
$ws.Range("A62").Value = "Monday, 21 July, 2025 1:44 PM"
